# WS_holdings.xlsx update:
#  - bump the "as of" date in the confidential disclaimer note (A16) from
#    2021-07-13 to 2021-07-14
#  - refresh the Weight / Percent Change figures (columns D & E, rows 2-13)
#    with the latest model-holdings snapshot

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet ships protected (no password prompt needed for this workbook's
# protection, but Excel still blocks writes until it is lifted). Unprotect,
# make the edits, then restore protection so the sheet stays locked as it
# was before.
$ws.Unprotect()

$ws.Range("A16").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-07-14 for illustrative purposes only and are subject to change."
# Setting a value with an embedded line break makes Excel auto-expand the row
# height; auto-fit it back down so row 16 keeps its original (default) height.
$ws.Rows(16).EntireRow.AutoFit()

$ws.Range("D2").Value = 0.02726201768194922
$ws.Range("E2").Value = 0.00705354280218029

$ws.Range("D3").Value = 0.02127037417605394
$ws.Range("E3").Value = 0.007444764649375823

$ws.Range("D4").Value = 0.05645843842897547
$ws.Range("E4").Value = 0.007292401787814429

$ws.Range("D5").Value = 0.1392631718988593
$ws.Range("E5").Value = -0.005265603221310089

$ws.Range("D6").Value = 0.02017064429550314
$ws.Range("E6").Value = -0.03140265177948354

$ws.Range("D7").Value = 0.1286098224676358
$ws.Range("E7").Value = -0.003654080389768555

$ws.Range("D8").Value = 0.08844879513285525
$ws.Range("E8").Value = -0.0024025133986324

$ws.Range("D9").Value = 0.02931765438786319
$ws.Range("E9").Value = -0.001517779705117128

$ws.Range("D10").Value = 0.103043827096309
$ws.Range("E10").Value = -0.004600345025876917

$ws.Range("D11").Value = 0.2968237694217881
$ws.Range("E11").Value = 0.004564315352697079

$ws.Range("D12").Value = 0.08933148501220735
$ws.Range("E12").Value = -0.00402561756633113

$ws.Range("D13").Value = 0.9999999999999998
$ws.Range("E13").Value = -0.0008101540212109848

$ws.Protect()
